$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J17").Value = 1087.2727
$ws.Range("H17").Value = 1087.2727
$ws.Range("L17").Value = 3261.8181
$ws.Range("N17").Value = -3597.8181
$ws.Range("J121").Value = 1413
$ws.Range("N121").Value = -7733
$ws.Range("L121").Value = 4239
$ws.Range("H121").Value = 1423.875
$ws.Range("H132").Value = 8136270
$ws.Range("M132").Value = -26326982
$ws.Range("I132").Value = 8776504
$ws.Range("K132").Value = 26329512
$ws.Range("J137").Value = 1746.5238
$ws.Range("L137").Value = 5239.5714
$ws.Range("I137").Value = 911.7027
$ws.Range("H137").Value = 1213.9656
$ws.Range("K137").Value = 2735.1081
$ws.Range("M137").Value = -185.1081000000004
$ws.Range("N137").Value = -10339.5714
$ws.Range("N138").Value = -15591.1145
$ws.Range("H138").Value = 1066.605
$ws.Range("L138").Value = 5311.1145
$ws.Range("J138").Value = 1770.3715
$ws.Range("I141").Value = 669.9
$ws.Range("K141").Value = 2009.7
$ws.Range("J141").Value = 2250
$ws.Range("N141").Value = -17110
$ws.Range("L141").Value = 6750
$ws.Range("M141").Value = 3170.3
$ws.Range("H141").Value = 933.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4090.48
$ws.Range("I32").Value = 3702.6365
$ws.Range("K32").Value = 3702.6365
$ws.Range("M32").Value = -3415.6365
$ws.Range("J74").Value = 1841
$ws.Range("N74").Value = -3589
$ws.Range("M74").Value = 142.64105
$ws.Range("L74").Value = 1841
$ws.Range("H74").Value = 939.4167
$ws.Range("K74").Value = 731.35895
$ws.Range("I74").Value = 731.35895
$ws.Range("J77").Value = 1841
$ws.Range("L77").Value = 9205
$ws.Range("N77").Value = -17941
$ws.Range("I77").Value = 731.35895
$ws.Range("H77").Value = 939.4167
$ws.Range("M77").Value = 711.20525
$ws.Range("K77").Value = 3656.79475
$ws.Range("M110").Value = 425
$ws.Range("H110").Value = 2188.8
$ws.Range("I110").Value = 1620
$ws.Range("K110").Value = 1620
$ws.Range("H132").Value = 2038.1471
$ws.Range("L132").Value = 5250
$ws.Range("J132").Value = 1750
$ws.Range("M132").Value = -3850.4228
$ws.Range("N132").Value = -10310
$ws.Range("I132").Value = 2126.8076
$ws.Range("K132").Value = 6380.4228
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("H22").Value = 1000
$ws.Range("M22").Value = -827
$ws.Range("K107").Value = 1103
$ws.Range("J107").Value = 1515.625
$ws.Range("M107").Value = 817
$ws.Range("I107").Value = 1103
$ws.Range("H107").Value = 1276.7368
$ws.Range("N107").Value = -5355.625
$ws.Range("L107").Value = 1515.625
$ws.Range("J134").Value = 16861.143
$ws.Range("L134").Value = 50583.429
$ws.Range("I134").Value = 1094.5667
$ws.Range("M134").Value = -748.7001
$ws.Range("H134").Value = 4077.4324
$ws.Range("N134").Value = -55653.429
$ws.Range("K134").Value = 3283.7001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K31").Value = 2314.9443
$ws.Range("N31").Value = -2570
$ws.Range("J31").Value = 1980
$ws.Range("H31").Value = 2242.1304
$ws.Range("I31").Value = 2314.9443
$ws.Range("L31").Value = 1980
$ws.Range("M31").Value = -2019.9443
$ws.Range("I34").Value = 2314.9443
$ws.Range("K34").Value = 2314.9443
$ws.Range("L34").Value = 1980
$ws.Range("H34").Value = 2242.1304
$ws.Range("J34").Value = 1980
$ws.Range("M34").Value = -2112.9443
$ws.Range("N34").Value = -2384
$ws.Range("L112").Value = 36456.223
$ws.Range("H112").Value = 34310.6
$ws.Range("J112").Value = 36456.223
$ws.Range("N112").Value = -39410.223
$ws.Range("H132").Value = 1721.55
$ws.Range("L132").Value = 8363.499899999999
$ws.Range("J132").Value = 2787.8333
$ws.Range("M132").Value = -2070.1469
$ws.Range("N132").Value = -13423.4999
$ws.Range("I132").Value = 1533.3823
$ws.Range("K132").Value = 4600.1469
$ws.Range("I134").Value = 893.0741
$ws.Range("M134").Value = -144.2223000000004
$ws.Range("H134").Value = 16130004
$ws.Range("K134").Value = 2679.2223
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L81").Value = 8974.940999999999
$ws.Range("J81").Value = 2991.647
$ws.Range("H81").Value = 2626.5
$ws.Range("N81").Value = -11220.941
$ws.Range("L84").Value = 26924.823
$ws.Range("N84").Value = -38156.823
$ws.Range("H84").Value = 2626.5
$ws.Range("J84").Value = 2991.647
$ws.Range("L122").Value = 7713
$ws.Range("N122").Value = -12613
$ws.Range("H122").Value = 759.7586
$ws.Range("K122").Value = 5598
$ws.Range("J122").Value = 857
$ws.Range("M122").Value = -3148
$ws.Range("I122").Value = 622
$ws.Range("I131").Value = 142857650
$ws.Range("N131").Value = -13938.0579
$ws.Range("K131").Value = 428572950
$ws.Range("L131").Value = 3858.0579
$ws.Range("M131").Value = -428567910
$ws.Range("J131").Value = 1286.0193
$ws.Range("H131").Value = 16950346
$ws.Range("I134").Value = 1539.5834
$ws.Range("M134").Value = 451.2497999999996
$ws.Range("H134").Value = 3535.3447
$ws.Range("K134").Value = 4618.7502
$ws.Range("L139").Value = 5023.928400000001
$ws.Range("N139").Value = -15303.9284
$ws.Range("H139").Value = 1771.3684
$ws.Range("J139").Value = 1674.6428
$ws.Range("M139").Value = -343.3747999999996
$ws.Range("K139").Value = 5483.3748
$ws.Range("I139").Value = 1827.7916
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("J70").Value = 33336726
$ws.Range("I70").Value = 31253232
$ws.Range("L70").Value = 33336726
$ws.Range("N70").Value = -33337266
$ws.Range("M70").Value = -31252962
$ws.Range("K70").Value = 31253232
$ws.Range("H70").Value = 32146158
$ws.Range("J73").Value = 33336726
$ws.Range("N73").Value = -33338598
$ws.Range("K73").Value = 31253232
$ws.Range("I73").Value = 31253232
$ws.Range("H73").Value = 32146158
$ws.Range("L73").Value = 33336726
$ws.Range("M73").Value = -31252296
$ws.Range("L106").Value = 65000
$ws.Range("H106").Value = 65000
$ws.Range("N106").Value = -67524
$ws.Range("J106").Value = 65000
$ws.Range("J113").Value = 1562.6
$ws.Range("H113").Value = 1489.7858
$ws.Range("N113").Value = -5902.6
$ws.Range("L113").Value = 1562.6
$ws.Range("L122").Value = 10698.9999
$ws.Range("N122").Value = -15598.9999
$ws.Range("H122").Value = 1908.55
$ws.Range("J122").Value = 3566.3333
$ws.Range("H132").Value = 2265.0952
$ws.Range("L132").Value = 10616.0001
$ws.Range("J132").Value = 3538.6667
$ws.Range("M132").Value = -2736.9998
$ws.Range("N132").Value = -15676.0001
$ws.Range("I132").Value = 1755.6666
$ws.Range("K132").Value = 5266.9998
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 934.2917
$ws.Range("N16").Value = -867.4286
$ws.Range("M16").Value = -931.8235
$ws.Range("K16").Value = 1101.8235
$ws.Range("J16").Value = 527.4286
$ws.Range("I16").Value = 1101.8235
$ws.Range("L16").Value = 527.4286
$ws.Range("J68").Value = 2200
$ws.Range("I68").Value = 2063.75
$ws.Range("M68").Value = -1314.75
$ws.Range("N68").Value = -3698
$ws.Range("K68").Value = 2063.75
$ws.Range("H68").Value = 2100.9092
$ws.Range("L68").Value = 2200
$ws.Range("H71").Value = 2100.9092
$ws.Range("K71").Value = 10318.75
$ws.Range("J71").Value = 2200
$ws.Range("M71").Value = -6574.75
$ws.Range("I71").Value = 2063.75
$ws.Range("L71").Value = 11000
$ws.Range("N71").Value = -18488
$ws.Range("K93").Value = 862.6667
$ws.Range("M93").Value = 385.3333
$ws.Range("H93").Value = 862.6667
$ws.Range("I93").Value = 862.6667
$ws.Range("K100").Value = 1074.5
$ws.Range("J100").Value = 2224.75
$ws.Range("M100").Value = -533.5
$ws.Range("L100").Value = 2224.75
$ws.Range("H100").Value = 1649.625
$ws.Range("I100").Value = 1074.5
$ws.Range("N100").Value = -3306.75
$ws.Range("H132").Value = 24754.047
$ws.Range("L132").Value = 256551
$ws.Range("J132").Value = 85517
$ws.Range("M132").Value = -1168.7096
$ws.Range("N132").Value = -261611
$ws.Range("I132").Value = 1232.9032
$ws.Range("K132").Value = 3698.7096
$ws.Range("K136").Value = 2344.61115
$ws.Range("M136").Value = 205.3888499999998
$ws.Range("H136").Value = 877.2456
$ws.Range("I136").Value = 781.53705
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I14").Value = 140161
$ws.Range("K14").Value = 140161
$ws.Range("J14").Value = 8571.143
$ws.Range("L14").Value = 8571.143
$ws.Range("H14").Value = 43200.05
$ws.Range("N14").Value = -8907.143
$ws.Range("M14").Value = -139993
$ws.Range("L126").Value = 6470.1432
$ws.Range("H126").Value = 37038510
$ws.Range("J126").Value = 2156.7144
$ws.Range("N126").Value = -11410.1432
$ws.Range("H132").Value = 4560.8423
$ws.Range("L132").Value = 4452.4998
$ws.Range("J132").Value = 1484.1666
$ws.Range("M132").Value = -15412.538
$ws.Range("N132").Value = -9512.4998
$ws.Range("I132").Value = 5980.846
$ws.Range("K132").Value = 17942.538
$ws.Range("J136").Value = 921
$ws.Range("L136").Value = 2763
$ws.Range("K136").Value = 1154.8929
$ws.Range("N136").Value = -7863
$ws.Range("M136").Value = 1395.1071
$ws.Range("H136").Value = 466.18182
$ws.Range("I136").Value = 384.9643
$ws.Range("J140").Value = 32374.8
$ws.Range("H140").Value = 32374.8
$ws.Range("N140").Value = -42734.8
$ws.Range("L140").Value = 32374.8
